$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 305.33334
$ws.Range("I28").Value = 239.5
$ws.Range("J28").Value = 387.625
$ws.Range("K28").Value = 239.5
$ws.Range("L28").Value = 387.625
$ws.Range("M28").Value = 245.5
$ws.Range("N28").Value = -1357.625

$ws.Range("H33").Value = 753.1111
$ws.Range("I33").Value = 504.25
$ws.Range("K33").Value = 504.25
$ws.Range("M33").Value = -275.25

$ws.Range("H40").Value = 1543.359
$ws.Range("I40").Value = 1136.7037
$ws.Range("J40").Value = 2458.3333
$ws.Range("K40").Value = 1136.7037
$ws.Range("L40").Value = 2458.3333
$ws.Range("M40").Value = -961.7037
$ws.Range("N40").Value = -2808.3333

$ws.Range("H64").Value = 3480.6667
$ws.Range("I64").Value = 3470
$ws.Range("J64").Value = 3502
$ws.Range("K64").Value = 3470
$ws.Range("L64").Value = 3502
$ws.Range("M64").Value = -3222
$ws.Range("N64").Value = -3998

$ws.Range("H67").Value = 3480.6667
$ws.Range("I67").Value = 3470
$ws.Range("J67").Value = 3502
$ws.Range("K67").Value = 3470
$ws.Range("L67").Value = 3502
$ws.Range("M67").Value = -2612
$ws.Range("N67").Value = -5218

$ws.Range("H74").Value = 5380.6
$ws.Range("I74").Value = 5634.3335
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 5634.3335
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -4698.3335
$ws.Range("N74").Value = -6872

$ws.Range("H77").Value = 5380.6
$ws.Range("I77").Value = 5634.3335
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 28171.6675
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -23491.6675
$ws.Range("N77").Value = -34360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 62502876
$ws.Range("I2").Value = 62502876
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 62502876
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -62502763
$ws.Range("N2").ClearContents()

$ws.Range("H116").Value = 62502876
$ws.Range("I116").Value = 62502876
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 62502876
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -62500582
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 417622.5
$ws.Range("I122").Value = 1792.5714
$ws.Range("J122").Value = 833452.4399999999
$ws.Range("K122").Value = 5377.7142
$ws.Range("L122").Value = 2500357.32
$ws.Range("M122").Value = -2927.7142
$ws.Range("N122").Value = -2505257.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 62502876
$ws.Range("I3").Value = 62502876
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 62502876
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -62502762
$ws.Range("N3").ClearContents()

$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -41694

$ws.Range("H64").Value = 622.9167
$ws.Range("I64").Value = 431.75
$ws.Range("J64").Value = 718.5
$ws.Range("K64").Value = 431.75
$ws.Range("L64").Value = 718.5
$ws.Range("M64").Value = -206.75
$ws.Range("N64").Value = -1168.5

$ws.Range("H67").Value = 622.9167
$ws.Range("I67").Value = 431.75
$ws.Range("J67").Value = 718.5
$ws.Range("K67").Value = 431.75
$ws.Range("L67").Value = 718.5
$ws.Range("M67").Value = 348.25
$ws.Range("N67").Value = -2278.5

$ws.Range("H105").Value = 2437.1667
$ws.Range("I105").Value = 2950
$ws.Range("J105").Value = 2180.75
$ws.Range("K105").Value = 2950
$ws.Range("L105").Value = 2180.75
$ws.Range("M105").Value = -1203
$ws.Range("N105").Value = -5674.75

$ws.Range("H134").Value = 3755.182
$ws.Range("I134").Value = 3412.0386
$ws.Range("J134").Value = 5029.7144
$ws.Range("K134").Value = 10236.1158
$ws.Range("L134").Value = 15089.1432
$ws.Range("M134").Value = -7701.1158
$ws.Range("N134").Value = -20159.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1668539.6
$ws.Range("I31").Value = 1787006.6
$ws.Range("K31").Value = 1787006.6
$ws.Range("M31").Value = -1786711.6

$ws.Range("H34").Value = 1668539.6
$ws.Range("I34").Value = 1787006.6
$ws.Range("K34").Value = 1787006.6
$ws.Range("M34").Value = -1786804.6

$ws.Range("H99").Value = 3800
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 4666.6665
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 4666.6665
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -7662.6665

$ws.Range("H126").Value = 3800
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -18939.9995

$ws.Range("H140").Value = 39700
$ws.Range("J140").Value = 39700
$ws.Range("L140").Value = 39700
$ws.Range("N140").Value = -50060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 314
$ws.Range("I8").Value = 314
$ws.Range("K8").Value = 942
$ws.Range("M8").Value = -803

$ws.Range("H103").Value = 2430.2727
$ws.Range("I103").Value = 2320
$ws.Range("J103").Value = 2522.1667
$ws.Range("K103").Value = 6960
$ws.Range("L103").Value = 7566.500100000001
$ws.Range("M103").Value = -6081
$ws.Range("N103").Value = -9324.500100000001

$ws.Range("H129").Value = 24507.521
$ws.Range("I129").Value = 4851.6665
$ws.Range("J129").Value = 31444.883
$ws.Range("K129").Value = 14554.9995
$ws.Range("L129").Value = 94334.649
$ws.Range("M129").Value = -9554.999500000002
$ws.Range("N129").Value = -104334.649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3993.8667
$ws.Range("I70").Value = 3998.6667
$ws.Range("K70").Value = 3998.6667
$ws.Range("M70").Value = -3728.6667

$ws.Range("H73").Value = 3993.8667
$ws.Range("I73").Value = 3998.6667
$ws.Range("K73").Value = 3998.6667
$ws.Range("M73").Value = -3062.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2248.6428
$ws.Range("I7").Value = 1700.125
$ws.Range("J7").Value = 2980
$ws.Range("K7").Value = 1700.125
$ws.Range("L7").Value = 2980
$ws.Range("M7").Value = -1588.125
$ws.Range("N7").Value = -3204

$ws.Range("H126").Value = 2248.6428
$ws.Range("I126").Value = 1700.125
$ws.Range("J126").Value = 2980
$ws.Range("K126").Value = 5100.375
$ws.Range("L126").Value = 8940
$ws.Range("M126").Value = -2630.375
$ws.Range("N126").Value = -13880
